# "changed cabbage to be a green"
#
# In the Ingredients sheet, the "Cabbage" row (row 28, category "Vegetable")
# is moved up to become the first row of the block (row 9, right after the
# "Greens" section), shifting the rows that used to be 9-27 down to 10-28,
# and its Category (column H) is changed from "Vegetable" to "Greens" so it
# now sits with the rest of the greens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The contiguous block that moves: rows 9 (Spinach) through 28 (Cabbage),
# columns A (Name) through H (Category).
$range = $ws.Range("A9:H28")
$data = $range.Value2

$firstRow = 1
$lastRow = 20
$lastCol = 8

# Stash the last row (old row 28 == Cabbage) before it gets overwritten.
$cabbage = New-Object 'object[]' $lastCol
for ($c = 1; $c -le $lastCol; $c++) {
    $cabbage[$c - 1] = $data[$lastRow, $c]
}

# Shift rows down by one: old row r (r = 19..1) -> new row r+1.
for ($r = $lastRow - 1; $r -ge $firstRow; $r--) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $data[$r + 1, $c] = $data[$r, $c]
    }
}

# Cabbage becomes the new first row of the block.
for ($c = 1; $c -le $lastCol; $c++) {
    $data[$firstRow, $c] = $cabbage[$c - 1]
}

# Recategorize cabbage as a green (column H).
$data[$firstRow, 8] = "Greens"

$range.Value2 = $data
